$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D").Insert()
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D7").Value2 = 43465
Write-Host "D7 Value2:" $ws.Range("D7").Value2
